$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.057.96"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").Value = "3.114.00"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.13"
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "625.60"
$ws.Range("E6").Value = "  -1.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.15"
$ws.Range("E7").Value = "  +6.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.375"
$ws.Range("E8").Value = "  +1.86%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.773"
$ws.Range("E10").Value = "  +6.71%  "

$ws.Range("B11").Value = "LidoStakedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D11").Value = "2.728.97"
$ws.Range("E11").Value = "  -13.27%  "

$ws.Range("E12").Value = "  +3.44%  "

$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.51"
$ws.Range("E14").Value = "  -2.48%  "

$ws.Range("D15").Value = "91.810.63"
$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("E16").Value = "  -2.02%  "

$ws.Range("D17").Value = "3.700.40"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").Value = "3.094.47"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.73"
$ws.Range("E19").Value = "  -0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000219"
$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.67"
$ws.Range("E21").Value = "  +1.81%  "

$ws.Range("E22").Value = "  +3.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "449.48"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.12"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "91.99"
$ws.Range("E25").Value = "  +1.18%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.63"
$ws.Range("E26").Value = "  -6.00%  "

$ws.Range("E27").Value = "  -4.36%  "

$ws.Range("D28").Value = "3.269.77"
$ws.Range("E28").Value = "  -0.98%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.185"
$ws.Range("E30").Value = "  +15.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.241"
$ws.Range("E31").Value = "  +18.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.35"
$ws.Range("E32").Value = "  -3.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.116"
$ws.Range("E33").Value = "  +34.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +1.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.169"
$ws.Range("E35").Value = "  +12.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.79"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.63"
$ws.Range("E37").Value = "  +6.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.06"
$ws.Range("E38").Value = "  +20.71%  "

$ws.Range("E39").Value = "  -1.44%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "492.00"
$ws.Range("E40").Value = "  -4.87%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("E41").Value = "  -5.92%  "

$ws.Range("E42").Value = "  -1.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.424"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.18"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.92"
$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "155.22"
$ws.Range("E47").Value = "  +2.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.695"
$ws.Range("E48").Value = "  -1.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.55"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("E50").Value = "  -0.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.72"
$ws.Range("E51").Value = "  -2.40%  "
